$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "recovery compression pants"
$ws.Range("A2").Value = "knee volleyball"
$ws.Range("A3").Value = "leggings black"
$ws.Range("A4").Value = "mens volleyball pads"
$ws.Range("A5").Value = "knee arthritis compression"
$ws.Range("A6").Value = "youth small baseball pants"
$ws.Range("A7").Value = "mens black pants"
$ws.Range("A8").Value = "knee high leggings"
$ws.Range("A9").Value = "knee pads arthritis"
$ws.Range("A10").Value = "basketball knee pads nike pro"
$ws.Range("A11").Value = "knee compression tights"
$ws.Range("A12").Value = "knee pad sleeve youth"
$ws.Range("A13").Value = "thermal compression pants men winter"
$ws.Range("A14").Value = "honey comb knee pad"
$ws.Range("A15").Value = "combat pants with knee pads"
$ws.Range("A16").Value = "deportivo para hombres"
$ws.Range("A17").Value = "compression spats"
$ws.Range("A18").Value = "red mens compression pants"
$ws.Range("A19").Value = "green basketball knee pads"
$ws.Range("A20").Value = "red basketball leggings"
$ws.Range("A21").Value = "underarmour mens tights"
$ws.Range("A22").Value = "three quarter compression pants men"
$ws.Range("A23").Value = "base layer pants men cold weather"
$ws.Range("A24").Value = "tesla base layer pants"
$ws.Range("A25").Value = "winter gym leggings"
$ws.Range("A26").Value = "thick spandex pants men"
$ws.Range("A27").Value = "knee pad leggings for basketball"
$ws.Range("A28").Value = "youth knee pads for basketball"
$ws.Range("A29").Value = "girls basketball knee pad leggings"
$ws.Range("A30").Value = "knee pads compression"
$ws.Range("A31").Value = "knee pad tights basketball"
$ws.Range("A32").Value = "weightlifting knee pads"
$ws.Range("A33").Value = "compression leggings youth"
$ws.Range("A34").Value = "hex pad knee"
$ws.Range("A35").Value = "knee pad soccer"
$ws.Range("A36").Value = "knee pad for running"
$ws.Range("A37").Value = "basketballs black"
$ws.Range("A38").Value = "boys soccer tights"
$ws.Range("A39").Value = "compression men leggings"
$ws.Range("A40").Value = "youth compression knee pads"
$ws.Range("A41").Value = "leggings with knee"
$ws.Range("A42").Value = "running compression pants men"
$ws.Range("A43").Value = "knee protector volleyball"
$ws.Range("A44").Value = "volleyball knee pads for men"
$ws.Range("A45").Value = "knee guard basketball"
$ws.Range("A46").Value = "black baseball pants boys"
$ws.Range("A47").Value = "compression pants boys"
$ws.Range("A48").Value = "knee pad adult"
$ws.Range("A49").Value = "little boy leggings"
$ws.Range("A50").Value = "capri legging pack"
$ws.Range("A51").Value = "men knee pads"
$ws.Range("A52").Value = "basketball compression"
$ws.Range("A53").Value = "boys compression tight"
$ws.Range("A54").Value = "baseball pants knee high"
$ws.Range("A55").Value = "small volleyball knee pads"
$ws.Range("A56").Value = "knee pads volleyball"
$ws.Range("A57").Value = "youth volleyball"
$ws.Range("A58").Value = "knee protector football"
$ws.Range("A59").Value = "knee pads protector"
$ws.Range("A60").Value = "snowboarding hip pads"
$ws.Range("A61").Value = "knee pad running"
$ws.Range("A62").Value = "men pads"
$ws.Range("A63").Value = "cycling compression tights"
$ws.Range("A64").Value = "capri leggings pack"
$ws.Range("A65").Value = "pants soccer"
$ws.Range("A66").Value = "running compression leggings"
$ws.Range("A67").Value = "boys black baseball pants"
$ws.Range("A68").Value = "leggings for boys"
$ws.Range("A69").Value = "youth girls compression pants"
$ws.Range("A70").Value = "athletic pants soccer"
$ws.Range("A71").Value = "knee guards basketball"
$ws.Range("A72").Value = "compression pants for hockey"
$ws.Range("A73").Value = "wrestling gear for men"
$ws.Range("A74").Value = "adult medium baseball pants"
$ws.Range("A75").Value = "knee pads for arthritis"
$ws.Range("A76").Value = "baseball boy pants"
$ws.Range("A77").Value = "knee protector pad"
$ws.Range("A78").Value = "knee protector soccer"
$ws.Range("A79").Value = "men sports compression"
$ws.Range("A80").Value = "basketball pants boys"
$ws.Range("A81").Value = "baseball pants"
$ws.Range("A82").Value = "knee protector pads"
$ws.Range("A83").Value = "knee pad for sports"
$ws.Range("A84").Value = "sports knee protectors"
$ws.Range("A85").Value = "youth knee pads volleyball"
$ws.Range("A86").Value = "bjj pants"
$ws.Range("A87").Value = "knee pads cycling"
$ws.Range("A88").Value = "boys youth basketball"
$ws.Range("A89").Value = "mens athletic pants tall"
$ws.Range("A90").Value = "legging pack"
$ws.Range("A91").Value = "girl volleyball knee pads"
$ws.Range("A92").Value = "running tights"
$ws.Range("A93").Value = "youth baseball pants"
$ws.Range("A94").Value = "boys sports pants"
$ws.Range("A95").Value = "compression for men"
$ws.Range("A96").Value = "basketballs youth size"
$ws.Range("A97").Value = "knee pads xl"
$ws.Range("A98").Value = "compression knee men"
$ws.Range("A99").Value = "knee guards for adults"
$ws.Range("A100").Value = "football knee pad"
